# barn/urls.xlsx — "bunch of barn files and images"
# Sheet1 gains a new "\barn" section (15 new url rows) and becomes the active
# sheet; Sheet2 gains one more row ("beurs" / "s-miel.gif") and is no longer
# the active sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 : new section header at row 3 ("\barn"), replacing what used to
# be the "yfcf" row there; the old row3/row4 pair ("yfcf"/"audio") shifts
# down two rows to become row5/row6, so clear the now-stale B3 and A4:B4. ---
$ws1.Range("A3").Value = "\barn"
$ws1.Range("B3").ClearContents()
$ws1.Range("A4").ClearContents()
$ws1.Range("B4").ClearContents()

# --- Sheet1 : re-assert rows 5-14 (old rows 3-12, now shifted down by two) ---
$sheet1Shifted = @(
    @(5,  "yfcf",      "https://web.archive.org/web/19990220153326im_/http://dewey.rug.ac.be/barn/tex/yfcf.html"),
    @(6,  "audio",     "https://web.archive.org/web/19990222072856im_/http://dewey.rug.ac.be/barn/tex/audio.html"),
    @(7,  "video",     "https://web.archive.org/web/19990428123223im_/http://dewey.rug.ac.be/barn/tex/video.html"),
    @(8,  "perfo",     "https://web.archive.org/web/19990220082926im_/http://dewey.rug.ac.be/barn/tex/perfo.html"),
    @(9,  "lyric",     "https://web.archive.org/web/19990219231714im_/http://dewey.rug.ac.be/barn/tex/lyric.html"),
    @(10, "max",       "https://web.archive.org/web/19990220013212im_/http://dewey.rug.ac.be/barn/tex/max.html"),
    @(11, "work",      "https://web.archive.org/web/19990220133854im_/http://dewey.rug.ac.be/barn/tex/work.html"),
    @(12, "yfcfprom",  "https://web.archive.org/web/19990220162412im_/http://dewey.rug.ac.be/barn/tex/yfcfprom.html"),
    @(13, "news",      "https://web.archive.org/web/19990220051455im_/http://dewey.rug.ac.be/barn/tex/news.html"),
    @(14, "neckprom",  "https://web.archive.org/web/19990501223642im_/http://dewey.rug.ac.be/barn/tex/neckprom.html")
)
foreach ($row in $sheet1Shifted) {
    $ws1.Range("A$($row[0])").Value = $row[1]
    $ws1.Range("B$($row[0])").Value = $row[2]
}

# --- Sheet1 : brand-new rows 15-27, the new "\barn" url listing ---
$sheet1New = @(
    @(15, "gold",     "https://web.archive.org/web/19990501210627im_/http://dewey.rug.ac.be/barn/tex/gold.html"),
    @(16, "claemit",  "https://web.archive.org/web/19990428110805im_/http://dewey.rug.ac.be/barn/tex/claemit.html"),
    @(17, "theresa",  "https://web.archive.org/web/19990501233653im_/http://dewey.rug.ac.be/barn/tex/theresa.html"),
    @(18, "nuns",     "https://web.archive.org/web/19990501230557im_/http://dewey.rug.ac.be/barn/tex/nuns.html"),
    @(19, "fiat",     "https://web.archive.org/web/19990501204011im_/http://dewey.rug.ac.be/barn/tex/fiat.html"),
    @(20, "dildo",    "https://web.archive.org/web/19990501200351im_/http://dewey.rug.ac.be/barn/tex/dildo.html"),
    @(21, "meno",     "https://web.archive.org/web/19990501221035im_/http://dewey.rug.ac.be/barn/tex/meno.html"),
    @(22, "herma",    "https://web.archive.org/web/19990501213417im_/http://dewey.rug.ac.be/barn/tex/herma.html"),
    @(23, "doppler",  "https://web.archive.org/web/19990501202926im_/http://dewey.rug.ac.be/barn/tex/doppler.html"),
    @(24, "piet",     "https://web.archive.org/web/19990220101839im_/http://dewey.rug.ac.be/barn/tex/piet.html"),
    @(25, "cv",       "https://web.archive.org/web/19990428113621im_/http://dewey.rug.ac.be/barn/tex/cv.html"),
    @(26, "beurs",    "https://web.archive.org/web/19991104064655im_/http://dewey.rug.ac.be/barn/tex/beurs.html"),
    @(27, "sputnick", "https://web.archive.org/web/19991104064845im_/http://dewey.rug.ac.be/barn/tex/sputnick.html")
)
foreach ($row in $sheet1New) {
    $ws1.Range("A$($row[0])").Value = $row[1]
    $ws1.Range("B$($row[0])").Value = $row[2]
}

# --- Sheet2 : one new row of images ("beurs" / "s-miel.gif") ---
$ws2.Range("A5").Value = "beurs"
$ws2.Range("B5").Value = "s-miel.gif"

# --- Selections: Sheet2's cursor moves to B6 but Sheet1 becomes (and stays)
# the active/selected tab with its cursor on B27, matching the saved view. ---
$ws2.Range("B6").Select()
$ws1.Range("B27").Select()
